# BAB V revision:
#  1. Merge the split runs "model yang lebih mu" + "takhir seperti" into a
#     single run reading "model yang lebih mutakhir seperti".
#  2. Add a reviewer comment on the "Kesimpulan" intro sentence.

$word.UserName = "PAULUS CAESARIO DITO PUTRA HARTONO"
$word.UserInitials = "PC"

$d = $word.ActiveDocument

# --- 1. Fix the accidentally split word "mutakhir" -------------------------
$rFix = $d.Content
$rFix.Find.Execute(
    "model yang lebih mutakhir seperti",  # FindText (matches across the two split runs)
    $true,                                 # MatchCase
    $false,                                # MatchWholeWord
    $false,                                # MatchWildcards
    $false,                                # MatchSoundsLike
    $false,                                # MatchAllWordForms
    $true,                                 # Forward
    1,                                      # Wrap (wdFindContinue)
    $false,                                # Format
    "model yang lebih mutakhir seperti",  # ReplaceWith
    2                                       # Replace (wdReplaceAll)
) | Out-Null

# --- 2. Add the reviewer comment --------------------------------------------
$rComment = $d.Content
$rComment.Find.Execute(
    "Berdasarkan serangkaian proses pengujian dan analisis yang telah dilakukan pada bab sebelumnya, dapat ditarik beberapa kesimpulan sebagai berikut: ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null

$comment = $d.Comments.Add($rComment, "Apakah sudah menjawab rumusan masalah")
$comment.Initial = "PC"
